# Auto-generated data refresh for Pandaemonium Profits workbook.
# Updates marketboard-derived price/profit columns (H:N) across several
# sheets to reflect the latest scheduled price pull.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1764.4242
$ws.Range("I70").Value = 1645.85
$ws.Range("J70").Value = 1946.8462
$ws.Range("K70").Value = 4937.549999999999
$ws.Range("L70").Value = 5840.5386
$ws.Range("M70").Value = -4667.549999999999
$ws.Range("N70").Value = -6380.5386
# Row 73
$ws.Range("H73").Value = 1764.4242
$ws.Range("I73").Value = 1645.85
$ws.Range("J73").Value = 1946.8462
$ws.Range("K73").Value = 4937.549999999999
$ws.Range("L73").Value = 5840.5386
$ws.Range("M73").Value = -4001.549999999999
$ws.Range("N73").Value = -7712.5386
# Row 112
$ws.Range("H112").Value = 5304.8335
$ws.Range("J112").Value = 1579.1333
$ws.Range("L112").Value = 4737.3999
$ws.Range("N112").Value = -6953.3999
# Row 121
$ws.Range("H121").Value = 1150.75
$ws.Range("J121").Value = 1000.8182
$ws.Range("L121").Value = 3002.4546
$ws.Range("N121").Value = -6496.4546

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 59
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36608
# Row 61
$ws.Range("H61").Value = 15156740
$ws.Range("I61").Value = 20838324
$ws.Range("K61").Value = 20838324
$ws.Range("M61").Value = -20838112
# Row 122
$ws.Range("H122").Value = 62500850
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 125000000
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 375000000
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -375004900
# Row 136
$ws.Range("H136").Value = 15156740
$ws.Range("I136").Value = 20838324
$ws.Range("K136").Value = 62514972
$ws.Range("M136").Value = -62512422

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 746971.4
$ws.Range("I105").Value = 895580
$ws.Range("J105").Value = 3928.5715
$ws.Range("K105").Value = 895580
$ws.Range("L105").Value = 3928.5715
$ws.Range("M105").Value = -893833
$ws.Range("N105").Value = -7422.5715

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1565.3846
$ws.Range("I16").Value = 1058.7693
$ws.Range("J16").Value = 2072
$ws.Range("K16").Value = 1058.7693
$ws.Range("L16").Value = 2072
$ws.Range("M16").Value = -771.7692999999999
$ws.Range("N16").Value = -2646
# Row 29
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20586
# Row 31
$ws.Range("H31").Value = 689301.3
$ws.Range("I31").Value = 6628.1665
$ws.Range("J31").Value = 1030637.9
$ws.Range("K31").Value = 6628.1665
$ws.Range("L31").Value = 1030637.9
$ws.Range("M31").Value = -6333.1665
$ws.Range("N31").Value = -1031227.9
# Row 34
$ws.Range("H34").Value = 689301.3
$ws.Range("I34").Value = 6628.1665
$ws.Range("J34").Value = 1030637.9
$ws.Range("K34").Value = 6628.1665
$ws.Range("L34").Value = 1030637.9
$ws.Range("M34").Value = -6426.1665
$ws.Range("N34").Value = -1031041.9
# Row 58
$ws.Range("H58").Value = 2335553.5
$ws.Range("I58").Value = 3368524
$ws.Range("J58").Value = 11369.5
$ws.Range("K58").Value = 3368524
$ws.Range("L58").Value = 11369.5
$ws.Range("M58").Value = -3368321
$ws.Range("N58").Value = -11775.5
# Row 113
$ws.Range("H113").Value = 1565.3846
$ws.Range("I113").Value = 1058.7693
$ws.Range("J113").Value = 2072
$ws.Range("K113").Value = 1058.7693
$ws.Range("L113").Value = 2072
$ws.Range("M113").Value = 1111.2307
$ws.Range("N113").Value = -6412
# Row 136
$ws.Range("H136").Value = 2335553.5
$ws.Range("I136").Value = 3368524
$ws.Range("J136").Value = 11369.5
$ws.Range("K136").Value = 10105572
$ws.Range("L136").Value = 34108.5
$ws.Range("M136").Value = -10103022
$ws.Range("N136").Value = -39208.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 3571.4285
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 3571.4285
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10714.2855
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -11052.2855
# Row 30
$ws.Range("H30").Value = 3571.4285
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 3571.4285
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 10714.2855
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -10918.2855
# Row 36
$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 500
$ws.Range("K36").Value = 1500
$ws.Range("M36").Value = -1331
# Row 96
$ws.Range("H96").Value = 3809.3333
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 4371.2
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 13113.6
$ws.Range("M96").Value = -941
$ws.Range("N96").Value = -17231.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 766.6667
$ws.Range("I12").Value = 650
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 650
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -480
$ws.Range("N12").Value = -1340
# Row 19
$ws.Range("H19").Value = 7364.857
$ws.Range("J19").Value = 12500.25
$ws.Range("L19").Value = 12500.25
$ws.Range("N19").Value = -12840.25
# Row 22
$ws.Range("H22").Value = 1333.1666
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 1674.75
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 1674.75
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -2264.75
# Row 27
$ws.Range("H27").Value = 1333.1666
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 1674.75
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 1674.75
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -1888.75
# Row 41
$ws.Range("H41").Value = 6025.7144
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 8036
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 8036
$ws.Range("M41").Value = -562
$ws.Range("N41").Value = -8912
# Row 42
$ws.Range("H42").Value = 19899.5
$ws.Range("I42").Value = 9999
$ws.Range("J42").Value = 29800
$ws.Range("K42").Value = 9999
$ws.Range("L42").Value = 29800
$ws.Range("M42").Value = -9436
$ws.Range("N42").Value = -30926
# Row 49
$ws.Range("H49").Value = 19899.5
$ws.Range("I49").Value = 9999
$ws.Range("J49").Value = 29800
$ws.Range("K49").Value = 9999
$ws.Range("L49").Value = 29800
$ws.Range("M49").Value = -9852
$ws.Range("N49").Value = -30094
# Row 122
$ws.Range("H122").Value = 5556.375
$ws.Range("J122").Value = 4612.1665
$ws.Range("L122").Value = 13836.4995
$ws.Range("N122").Value = -18736.4995
# Row 136
$ws.Range("H136").Value = 3156.2927
$ws.Range("I136").Value = 2089.3157
$ws.Range("J136").Value = 4077.7727
$ws.Range("K136").Value = 6267.9471
$ws.Range("L136").Value = 12233.3181
$ws.Range("M136").Value = -3717.9471
$ws.Range("N136").Value = -17333.3181

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 26083.334
$ws.Range("J20").Value = 26083.334
$ws.Range("L20").Value = 26083.334
$ws.Range("N20").Value = -26563.334
# Row 47
$ws.Range("H47").Value = 123666.664
$ws.Range("J47").Value = 123666.664
$ws.Range("L47").Value = 123666.664
$ws.Range("N47").Value = -124810.664
# Row 122
$ws.Range("H122").Value = 2162.375
$ws.Range("I122").Value = 2179.8
$ws.Range("J122").Value = 2133.3333
$ws.Range("K122").Value = 6539.400000000001
$ws.Range("L122").Value = 6399.999899999999
$ws.Range("M122").Value = -4089.400000000001
$ws.Range("N122").Value = -11299.9999

